# Update values in column B for specific rows per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = 8.5725
    7   = 5.102400000000002
    16  = 5.086399999999998
    28  = 6.168900000000001
    29  = 5.090300000000004
    32  = 6.861399999999999
    40  = 8.956399999999997
    52  = 5.480399999999999
    57  = 5.008899999999996
    66  = 5.750399999999996
    100 = 5.158200000000003
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
